$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.635.01'
$ws.Range('E2').Value = '  -0.67%  '
$ws.Range('D3').Value = '2.297.45'
$ws.Range('E3').Value = '  -0.39%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '301.18'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '95.69'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.41%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.508'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.69%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.493'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.62'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.25'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.91%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0786'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.41%  '
$ws.Range('E13').Value = '  -0.47%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.76'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.27%  '
$ws.Range('D15').Value = '2.653.03'
$ws.Range('E15').Value = '  -0.53%  '
$ws.Range('D16').Value = '2.305.21'
$ws.Range('E16').Value = '  -0.19%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.783'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('D18').Value = '42.567.88'
$ws.Range('E18').Value = '  -0.69%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.29'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.54%  '
$ws.Range('D20').Value = '0.0₃0890'
$ws.Range('E20').Value = '  -1.43%  '
$ws.Range('E21').Value = '  -0.93%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.83'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.26'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.81%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '234.85'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.82%  '
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('E26').Value = '  -2.65%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.53'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.68%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.37'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +14.34%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '164.32'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.71%  '
$ws.Range('E30').Value = '  -0.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '32.07'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.35%  '
$ws.Range('E32').Value = '  -0.07%  '
$ws.Range('E33').Value = '  -0.82%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.53'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.24%  '
$ws.Range('E35').Value = '  -6.95%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0701'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.87%  '
$ws.Range('E37').Value = '  -3.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.1000'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.74%  '
$ws.Range('E39').Value = '  -0.80%  '
$ws.Range('E40').Value = '  -1.42%  '
$ws.Range('E41').Value = '  -0.71%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '20.18'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +10.15%  '
$ws.Range('D43').Value = '1.966.95'
$ws.Range('E43').Value = '  -1.95%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.46'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.49%  '
$ws.Range('E45').Value = '  -0.65%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.03'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.86%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.75'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.08%  '
$ws.Range('B48').Value = 'HuobiToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.92'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.50%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '2.527.21'
$ws.Range('E49').Value = '  -0.30%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '53.13'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.37%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '71.27'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.87%  '
